$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 10 and row 11 (two species records that were
# listed in the wrong order). After the change, row 10 holds what used to be
# row 11's record (Garnlav / Alectoria sarmentosa) and row 11 holds what used
# to be row 10's record (Tretåig hackspett / Picoides tridactylus).

# ---- Row 10: becomes the former row 11 content ----
$ws.Range("A10").Value = 130752527
$ws.Range("B10").Value = 79243
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."

# These fields only existed on row 10 before the swap; row 11 never had
# values here, so they are cleared out on row 10 now.
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("AC10").Value = ""

$ws.Range("Q10").Value = 490673
$ws.Range("R10").Value = 6763435
$ws.Range("Z10").Value = "11:43"
$ws.Range("AB10").Value = "11:43"
$ws.Range("AW10").Value = "Håkan Thenander"
$ws.Range("AX10").Value = "Håkan Thenander, Bo karlstens"

# ---- Row 11: becomes the former row 10 content ----
$ws.Range("A11").Value = 130789462
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"

# These fields now move onto row 11, matching the former row 10 values.
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = "äldre spår"
$ws.Range("N11").Value = ""
$ws.Range("AC11").Value = "Äldre ring hack på tall"

$ws.Range("Q11").Value = 490705
$ws.Range("R11").Value = 6763439
$ws.Range("Z11").Value = "12:09"
$ws.Range("AB11").Value = "12:09"
$ws.Range("AW11").Value = "Bo karlstens"
$ws.Range("AX11").Value = "Bo karlstens, Håkan Thenander"
